# Revise config file handling
# Appends three new data rows (144-146) to each of the four worksheets
# (MID_LFT_#1, MID_LFT_#2, MID_PLT_#1, MID_PLT_#2), mirroring the existing
# per-sheet config-record pattern already present in rows 1-143.

$wb = $excel.ActiveWorkbook

$dateFmt = "YYYY-MM-DD HH:MM:SS"

function Set-ConfigRow {
    param(
        $ws,
        [int]$row,
        [double]$aVal,
        [string]$bVal,
        [string]$cVal,
        [string]$dVal,
        [string]$eVal,
        [double]$fVal,
        [double]$gVal,
        [double]$hVal,
        [double]$iVal
    )

    $ws.Cells.Item($row, 1).Value = $aVal
    $ws.Cells.Item($row, 1).NumberFormat = $dateFmt

    $ws.Cells.Item($row, 2).Value = $bVal
    $ws.Cells.Item($row, 3).Value = $cVal
    $ws.Cells.Item($row, 4).Value = $dVal
    $ws.Cells.Item($row, 5).Value = $eVal

    $ws.Cells.Item($row, 6).Value = $fVal
    $ws.Cells.Item($row, 7).Value = $gVal
    $ws.Cells.Item($row, 8).Value = $hVal
    $ws.Cells.Item($row, 9).Value = $iVal
}

# Excel's PI-style scientific-notation literal (e.g. 5.68...e+23) is parsed
# via an explicit [double] string cast below, since bare `e+23` exponent
# suffixes on a numeric literal aren't accepted by the script parser.
$g1 = [double]"5.68631262647113e+23"
$g2 = [double]"5.68432987514711e+23"

# ---- Sheet 1: MID_LFT_#1 ----
$ws1 = $wb.Worksheets.Item(1)

Set-ConfigRow $ws1 144 45930.46016203704 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c," "0x00,0xD4" "0x07" 400 $g1 220 7
Set-ConfigRow $ws1 145 45931.4603125     "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c," "0x00,0xD4" "0x07" 400 $g1 220 7
Set-ConfigRow $ws1 146 45932.46620370371 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c," "0x00,0xD4" "0x07" 400 $g1 216 7

# ---- Sheet 2: MID_LFT_#2 ----
$ws2 = $wb.Worksheets.Item(2)

Set-ConfigRow $ws2 144 45930.46016203704 "0x01,0x7c" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x00,0xF0" "0x19" 380 $g2 240 25
Set-ConfigRow $ws2 145 45931.4603125     "0x01,0x7c" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x00,0xF0" "0x19" 380 $g2 240 25
Set-ConfigRow $ws2 146 45932.46620370371 "0x01,0x7c" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x00,0xF0" "0x19" 380 $g2 240 25

# ---- Sheet 3: MID_PLT_#1 ----
$ws3 = $wb.Worksheets.Item(3)

Set-ConfigRow $ws3 144 45930.46016203704 "0x00,0x6e" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," "0x00,0x54" "0x15" 110 $g1 84 15
Set-ConfigRow $ws3 145 45931.4603125     "0x00,0x6e" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," "0x00,0x54" "0x15" 110 $g1 84 15
Set-ConfigRow $ws3 146 45932.46620370371 "0x00,0x6e" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," "0x00,0x53" "0x15" 110 $g1 83 15

# ---- Sheet 4: MID_PLT_#2 ----
$ws4 = $wb.Worksheets.Item(4)

Set-ConfigRow $ws4 144 45930.46016203704 "0x00,0x82" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," "0x00,0x6C" "0x9" 130 $g1 108 9
Set-ConfigRow $ws4 145 45931.4603125     "0x00,0x82" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," "0x00,0x6C" "0x9" 130 $g1 108 9
Set-ConfigRow $ws4 146 45932.46620370371 "0x00,0x82" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," "0x00,0x6B" "0x9" 130 $g1 107 9
